# Apply the edits described by the diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row text (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Fix capitalization of connecting "de" -> "De" in a handful of place names
$ws.Range("A6").Value = "Ciudad De México"
$ws.Range("B10").Value = "Atlamajalcingo Del Monte"
$ws.Range("B11").Value = "Mártir De Cuilapan"
$ws.Range("B15").Value = "Unión De Tula"
$ws.Range("B22").Value = "Jalpan De Serra"
$ws.Range("B31").Value = "Poza Rica De Hidalgo"

# 3. Remove the trailing footnote/metadata rows (37-41)
$ws.Range("A37:A41").EntireRow.Delete()
